$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The syllabus table gained a new "Programa resumido:" (summary
# syllabus) content row. That pushes every row from the old row 14
# ("Short syllabus:") down by one. Insert a blank row above the old
# row 14 to make room, then rewrite every affected cell explicitly so
# the final content/formatting matches the target state regardless of
# what Excel copied down during the insert.
# ------------------------------------------------------------------
$ws.Rows.Item(14).Insert()

# ---- Row 10: Objetivos: content updated ----
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "Avaliar casos de impacto ambiental que proporcionaram singularidades em algumas de suas etapas."
$ws.Range("C10").Value = "Avaliar casos de impacto ambiental que proporcionaram singularidades em algumas de suas etapas."

# ---- Row 11: Objectives: (label only) ----
$ws.Range("A11").Value = "Objectives:"

# ---- Row 12: Docentes responsáveis: (label only) ----
$ws.Range("A12").Value = "Docentes responsáveis:"

# ---- Row 13: professor name moved here, no label in column A ----
$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = "5840938 - Marcelo Rodrigues de Holanda"
$ws.Range("C13").Value = "5840938 - Marcelo Rodrigues de Holanda"

# ---- Row 14 (new): Programa resumido: + new summary text ----
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "Estudos de caso: avaliar a singularidade do caso apresentado, com a necessária identificação das etapas e das peculiaridades que os fizeram próprios a serem aplicados em estudo de caso."
$ws.Range("C14").Value = "Estudos de caso: avaliar a singularidade do caso apresentado, com a necessária identificação das etapas e das peculiaridades que os fizeram próprios a serem aplicados em estudo de caso."

# ---- Row 15: Short syllabus: (label only) ----
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()

# ---- Row 16: Programa: + new syllabus text ----
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "Estudos de casos específicos e as etapas necessárias na avaliação de um impacto ambiental."
$ws.Range("C16").Value = "Estudos de casos específicos e as etapas necessárias na avaliação de um impacto ambiental."

# ---- Row 17: Syllabus: (label only) ----
$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

# ---- Row 18: Avaliação: (label only) ----
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18").ClearContents()
$ws.Range("C18").ClearContents()

# ---- Row 19: Método: + teaching method text ----
$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "Aula expositiva e exercícios dirigidos."
$ws.Range("C19").Value = "Aula expositiva e exercícios dirigidos."

# ---- Row 20: Critério: + grading criteria text ----
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "Média ponderada de exercícios e provas."
$ws.Range("C20").Value = "Média ponderada de exercícios e provas."

# ---- Row 21: Norma de recuperação: + recovery rule text ----
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "Prova única com nota igual ou superior a 5,0."
$ws.Range("C21").Value = "Prova única com nota igual ou superior a 5,0."

# ---- Row 22: Bibliografia: + bibliography text ----
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "Estudos de caso: EPIA de origem."
$ws.Range("C22").Value = "Estudos de caso: EPIA de origem."

# ------------------------------------------------------------------
# Row heights: the table alternates between default height and
# explicit 60pt / 120pt rows depending on how much text each entry
# holds. Re-apply the exact target heights; AutoFit() clears the
# customHeight flag for rows that should go back to the default.
# ------------------------------------------------------------------
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
